# "update task assign and deadlines"
#
# Sheet2 column E holds task deadlines. Row 22's deadline ("8/4/2013 12:00AM")
# is corrected to "6/4/2013 12:00AM", and rows 23-48 (which shared that same
# string) are moved to a new deadline, "6/6/2013 12:00AM".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Row 22 keeps its own (corrected) date.
$ws.Cells.Item(22, 5).Value = "6/4/2013 12:00AM"

# Rows 23-48 move to the new shared deadline.
for ($r = 23; $r -le 48; $r++) {
    $ws.Cells.Item($r, 5).Value = "6/6/2013 12:00AM"
}

# Reflect the editor's final scroll/selection state: looking at the
# reassigned block, selected E23:E48 with E23 active, scrolled so row 27
# is at the top of the window.
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E23:E48").Select()
